$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Wednesday 20170111 trade-plan row, inserted into the previously-blank row 3.
$ws.Range("B3").Value = "Tuesday"
$ws.Range("C3").Value = 20170111
$ws.Range("D3").Value = "Wednesday"
$ws.Range("E3").Value = 'NG had a big drop on Monday but met strong support at 3.103, only posted one 1-hour green candle the whole day before 4pm. The weather forecast is still warm (but speculaters are saying Canada''s code air will travel south), with the huge UGAZ trading volume, the NG keep rising over the night trading hours, and rised around 6% on Tuesday and totally shadow Monday''s red candle. With this trend, the NG prices will very likely to keep going again tomorrow on Wednesday, especially expectin a big withdraw for Thursday''s report. Still not sure Tuesday''s surge was due to anticipation of weather chagne or more from the profit taking/short recovering angle, but it''s probably conservative to sell some UGAZ to lock some profit. Or wait for the price to top, sell and buy back. '
$ws.Range("F3").Value = 'Should be more patient on buying UGAZ while NG price is reducing during the course of day. It is good to get in when having a strong feeling the trend will probably reverse, telling by the UGAZ trading volume, but patience will pay off if wait longer, there will be better enter prices. Normally, a big bottom will give you enough time to enter, because it''s big! '
$ws.Range("G3").Value = 20170109

# D3 loses the inherited wrap+left-align formatting shared by the rest of the
# row (it ends up back at the plain wrap-only style once re-typed).
$ws.Range("D3").HorizontalAlignment = 1

# A3 is a brand new cell in this row; give it the same wrap + left-aligned
# look as the rest of the row before filling in its date value.
$ws.Range("A3").WrapText = $true
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").Value = 20170110

# Row grows tall enough to show the new multi-line notes.
$ws.Range("A3:G3").RowHeight = 100.8

$ws.Range("E4").Select()
